# Apply the routing/date/time edits described by the diff.
# NOTE: Find.Execute with Replace:=wdReplaceAll (2) on a whole-document
# range (e.g. $d.Content) replaces every match in the document at once.
# Several of the target strings ("Татаркин", "18.01.2024") occur twice
# in this document with different replacement values, so each edit is
# scoped to the specific paragraph that contains it.

$d = $word.ActiveDocument

# Paragraph 4: "от Татаркин" -> "от Карыгин"
$p = $d.Paragraphs.Item(4).Range
$p.Find.Execute("Татаркин", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Карыгин", 2)

# Paragraph 6: "Прошу предоставить мне 18.01.2024 c 08:00 до 09:00 в счет отпуска."
$p = $d.Paragraphs.Item(6).Range
$p.Find.Execute("18.01.2024", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "23.01.2024", 2)

$p = $d.Paragraphs.Item(6).Range
$p.Find.Execute("08:00", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "13:00", 2)

$p = $d.Paragraphs.Item(6).Range
$p.Find.Execute("09:00", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "17:00", 2)

# Paragraph 8: the standalone date line near the signature
$p = $d.Paragraphs.Item(8).Range
$p.Find.Execute("18.01.2024", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "22.01.2024", 2)

# Paragraph 12: signature line "Татаркин" -> "Карыгин"
$p = $d.Paragraphs.Item(12).Range
$p.Find.Execute("Татаркин", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Карыгин", 2)
